$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Add new row of data (row 9)
$ws.Range("A9").Value = "721. Accounts Merge"
$ws.Range("B9").Value = "Medium"
$ws.Range("C9").Value = "Union Find"
$ws.Range("E9").Value = "https://leetcode.com/problems/accounts-merge/solutions/1601980/java-solution-using-unionfind-beats-99-87-of-submissions/?envType=study-plan-v2&envId=graph-theory "
$ws.Range("D9").Value = "Union Find, use hash for unique emails, and a hash for unique accounts, then collect the emails from the accounts hash and add account name at index 0 for the result list."

# Add the hyperlink for E9
$ws.Hyperlinks.Add($ws.Range("E9"), "https://leetcode.com/problems/accounts-merge/solutions/1601980/java-solution-using-unionfind-beats-99-87-of-submissions/?envType=study-plan-v2&envId=graph-theory ") | Out-Null

# Copy style formatting from row 8 to row 9 (B column orange fill, E column green + hyperlink style)
$ws.Range("B8").Copy() | Out-Null
$ws.Range("B9").PasteSpecial(-4122) | Out-Null  # xlPasteFormats
$ws.Range("E8").Copy() | Out-Null
$ws.Range("E9").PasteSpecial(-4122) | Out-Null  # xlPasteFormats
$excel.CutCopyMode = 0

# Update the sheet view: remove topLeftCell override, update selection to D10
$ws.Application.ActiveWindow.ScrollColumn = 1
$ws.Range("D10").Select() | Out-Null
